# Apply updated "dSF" (column F) values as per repull/recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    6  = 0
    9  = 1
    12 = 2
    14 = 7
    15 = -2
    17 = -3
    18 = -1
    19 = 4
    20 = 2
    24 = 1
    27 = -3
    28 = 1
    37 = -1
    39 = -2
    40 = -3
    41 = -2
    43 = -3
    46 = -2
    51 = -8
    56 = -3
    58 = -3
    59 = -9
    60 = 2
    61 = -3
    62 = -4
    64 = -5
    65 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
